# Weekly update: insert the new daily price record for Mandarina (Murcott,
# Segunda) at row 41, pushing the existing historical rows (old 41-96) down
# by one to rows 42-97.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 41; Excel shifts rows 41:96 down to 42:97 and
# carries the row's number formatting (date style on column D) along.
$ws.Rows.Item(41).Insert()

# Populate the newly inserted row 41 with the new observation.
$ws.Range("A41").Value = 1
$ws.Range("B41").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C41").Value = "Arica y Parinacota"
$ws.Range("D41").Value = 44679
$ws.Range("E41").Value = 15
$ws.Range("F41").Value = "Fruta"
$ws.Range("G41").Value = 100102
$ws.Range("H41").Value = "Cítricos"
$ws.Range("I41").Value = 100102004
$ws.Range("J41").Value = "Mandarina"
$ws.Range("K41").Value = "Murcott"
$ws.Range("L41").Value = "Segunda"
$ws.Range("M41").Value = 250
$ws.Range("N41").Value = 19000
$ws.Range("O41").Value = 20000
$ws.Range("P41").Value = 19500
$ws.Range("Q41").Value = "$/caja 20 kilos"
$ws.Range("R41").Value = "Región de Coquimbo"
$ws.Range("S41").Value = 975
$ws.Range("T41").Value = 20
